$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 348.01
$ws.Range("I15").Value = 348.01
$ws.Range("K15").Value = 1044.03
$ws.Range("M15").Value = -875.03
# Row 33
$ws.Range("H33").Value = 547.7273
$ws.Range("I33").Value = 645.4666999999999
$ws.Range("J33").Value = 338.2857
$ws.Range("K33").Value = 645.4666999999999
$ws.Range("L33").Value = 338.2857
$ws.Range("M33").Value = -416.4666999999999
$ws.Range("N33").Value = -796.2857
# Row 125
$ws.Range("H125").Value = 2757.1428
$ws.Range("I125").Value = 1866.6666
$ws.Range("K125").Value = 16799.9994
$ws.Range("M125").Value = -14339.9994
# Row 137
$ws.Range("H137").Value = 6501.875
$ws.Range("I137").Value = 7374.4
$ws.Range("J137").Value = 3385.7144
$ws.Range("K137").Value = 22123.2
$ws.Range("L137").Value = 10157.1432
$ws.Range("M137").Value = -19573.2
$ws.Range("N137").Value = -15257.1432
# Row 138
$ws.Range("H138").Value = 2660.4358
$ws.Range("I138").Value = 1294.326
$ws.Range("J138").Value = 4624.2188
$ws.Range("K138").Value = 3882.978
$ws.Range("L138").Value = 13872.6564
$ws.Range("M138").Value = 1257.022
$ws.Range("N138").Value = -24152.6564
# Row 140
$ws.Range("H140").Value = 38890
$ws.Range("I140").Value = 29000
$ws.Range("J140").Value = 48780
$ws.Range("K140").Value = 29000
$ws.Range("L140").Value = 48780
$ws.Range("M140").Value = -23820
$ws.Range("N140").Value = -59140
# Row 141
$ws.Range("H141").Value = 1229.6
$ws.Range("I141").Value = 723.26086
$ws.Range("J141").Value = 7052.5
$ws.Range("K141").Value = 2169.78258
$ws.Range("L141").Value = 21157.5
$ws.Range("M141").Value = 3010.21742
$ws.Range("N141").Value = -31517.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1211754
$ws.Range("I32").Value = 1347690.5
$ws.Range("J32").Value = 7744.857
$ws.Range("K32").Value = 1347690.5
$ws.Range("L32").Value = 7744.857
$ws.Range("M32").Value = -1347403.5
$ws.Range("N32").Value = -8318.857
# Row 45
$ws.Range("H45").Value = 1158.8
$ws.Range("I45").Value = 1480
$ws.Range("J45").Value = 1078.5
$ws.Range("K45").Value = 1480
$ws.Range("L45").Value = 1078.5
$ws.Range("M45").Value = -1103
$ws.Range("N45").Value = -1832.5
# Row 122
$ws.Range("H122").Value = 3320.976
$ws.Range("I122").Value = 3426.2593
$ws.Range("J122").Value = 3131.4666
$ws.Range("K122").Value = 10278.7779
$ws.Range("L122").Value = 9394.399800000001
$ws.Range("M122").Value = -7828.777900000001
$ws.Range("N122").Value = -14294.3998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
# Row 99
$ws.Range("H99").Value = 3358.76
$ws.Range("I99").Value = 4223.0625
$ws.Range("J99").Value = 1822.2222
$ws.Range("K99").Value = 4223.0625
$ws.Range("L99").Value = 1822.2222
$ws.Range("M99").Value = -2725.0625
$ws.Range("N99").Value = -4818.2222
# Row 105
$ws.Range("H105").Value = 2101
$ws.Range("I105").Value = 2403.3333
$ws.Range("J105").Value = 1971.4286
$ws.Range("K105").Value = 2403.3333
$ws.Range("L105").Value = 1971.4286
$ws.Range("M105").Value = -656.3332999999998
$ws.Range("N105").Value = -5465.4286

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2445.635
$ws.Range("I31").Value = 1611.9791
$ws.Range("J31").Value = 5113.3335
$ws.Range("K31").Value = 1611.9791
$ws.Range("L31").Value = 5113.3335
$ws.Range("M31").Value = -1316.9791
$ws.Range("N31").Value = -5703.3335
# Row 34
$ws.Range("H34").Value = 2445.635
$ws.Range("I34").Value = 1611.9791
$ws.Range("J34").Value = 5113.3335
$ws.Range("K34").Value = 1611.9791
$ws.Range("L34").Value = 5113.3335
$ws.Range("M34").Value = -1409.9791
$ws.Range("N34").Value = -5517.3335
# Row 58
$ws.Range("H58").Value = 3490.9778
$ws.Range("I58").Value = 3775
$ws.Range("J58").Value = 2791.8462
$ws.Range("K58").Value = 3775
$ws.Range("L58").Value = 2791.8462
$ws.Range("M58").Value = -3572
$ws.Range("N58").Value = -3197.8462
# Row 132
$ws.Range("H132").Value = 2122.25
$ws.Range("I132").Value = 1197.7037
$ws.Range("J132").Value = 4042.4614
$ws.Range("K132").Value = 3593.1111
$ws.Range("L132").Value = 12127.3842
$ws.Range("M132").Value = -1063.1111
$ws.Range("N132").Value = -17187.3842
# Row 136
$ws.Range("H136").Value = 3490.9778
$ws.Range("I136").Value = 3775
$ws.Range("J136").Value = 2791.8462
$ws.Range("K136").Value = 11325
$ws.Range("L136").Value = 8375.5386
$ws.Range("M136").Value = -8775
$ws.Range("N136").Value = -13475.5386

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 74
$ws.Range("H74").Value = 2956.8572
$ws.Range("J74").Value = 3399.6667
$ws.Range("L74").Value = 10199.0001
$ws.Range("N74").Value = -12321.0001
# Row 77
$ws.Range("H77").Value = 2956.8572
$ws.Range("J77").Value = 3399.6667
$ws.Range("L77").Value = 30597.0003
$ws.Range("N77").Value = -41205.0003
# Row 120
$ws.Range("H120").Value = 9211
$ws.Range("I120").Value = 0
$ws.Range("J120").Value = 9211
$ws.Range("K120").Value = 0
$ws.Range("L120").Value = 27633
$ws.Range("M120").ClearContents()
$ws.Range("N120").Value = -37309
# Row 131
$ws.Range("H131").Value = 2549.415
$ws.Range("I131").Value = 2627.6924
$ws.Range("J131").Value = 2523.975
$ws.Range("K131").Value = 7883.0772
$ws.Range("L131").Value = 7571.924999999999
$ws.Range("M131").Value = -2843.0772
$ws.Range("N131").Value = -17651.925
# Row 132
$ws.Range("H132").Value = 12625.375
$ws.Range("I132").Value = 7800.6
$ws.Range("J132").Value = 20666.666
$ws.Range("K132").Value = 70205.40000000001
$ws.Range("L132").Value = 185999.994
$ws.Range("M132").Value = -67675.40000000001
$ws.Range("N132").Value = -191059.994

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4647.727
$ws.Range("I70").Value = 4444.4443
$ws.Range("J70").Value = 4788.4614
$ws.Range("K70").Value = 4444.4443
$ws.Range("L70").Value = 4788.4614
$ws.Range("M70").Value = -4174.4443
$ws.Range("N70").Value = -5328.4614
# Row 73
$ws.Range("H73").Value = 4647.727
$ws.Range("I73").Value = 4444.4443
$ws.Range("J73").Value = 4788.4614
$ws.Range("K73").Value = 4444.4443
$ws.Range("L73").Value = 4788.4614
$ws.Range("M73").Value = -3508.4443
$ws.Range("N73").Value = -6660.4614
# Row 80
$ws.Range("H80").Value = 11113.125
$ws.Range("I80").Value = 12315
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 12315
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -11317
$ws.Range("N80").Value = -4696
# Row 83
$ws.Range("H83").Value = 11113.125
$ws.Range("I83").Value = 12315
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 61575
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -56583
$ws.Range("N83").Value = -23484
# Row 102
$ws.Range("H102").Value = 3500.258
$ws.Range("I102").Value = 1826
$ws.Range("J102").Value = 6151.1665
$ws.Range("K102").Value = 1826
$ws.Range("L102").Value = 6151.1665
$ws.Range("M102").Value = -204
$ws.Range("N102").Value = -9395.166499999999
# Row 122
$ws.Range("H122").Value = 1619.7222
$ws.Range("I122").Value = 1497.1875
$ws.Range("K122").Value = 4491.5625
$ws.Range("M122").Value = -2041.5625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1076.4468
$ws.Range("I16").Value = 912.2195
$ws.Range("J16").Value = 2198.6667
$ws.Range("K16").Value = 912.2195
$ws.Range("L16").Value = 2198.6667
$ws.Range("M16").Value = -742.2195
$ws.Range("N16").Value = -2538.6667
# Row 40
$ws.Range("H40").Value = 2964.7307
$ws.Range("I40").Value = 2367.7896
$ws.Range("J40").Value = 4585
$ws.Range("K40").Value = 2367.7896
$ws.Range("L40").Value = 4585
$ws.Range("M40").Value = -2231.7896
$ws.Range("N40").Value = -4857
# Row 55
$ws.Range("H55").Value = 192.21053
$ws.Range("I55").Value = 166.15384
$ws.Range("J55").Value = 248.66667
$ws.Range("K55").Value = 166.15384
$ws.Range("L55").Value = 248.66667
$ws.Range("M55").Value = 6.846159999999998
$ws.Range("N55").Value = -594.6666700000001
# Row 122
$ws.Range("H122").Value = 2937.4666
$ws.Range("I122").Value = 2933.818
$ws.Range("J122").Value = 2947.5
$ws.Range("K122").Value = 8801.454000000002
$ws.Range("L122").Value = 8842.5
$ws.Range("M122").Value = -6351.454000000002
$ws.Range("N122").Value = -13742.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 673.1539
$ws.Range("I107").Value = 626.6667
$ws.Range("J107").Value = 777.75
$ws.Range("K107").Value = 1880.0001
$ws.Range("L107").Value = 2333.25
$ws.Range("M107").Value = 39.99990000000003
$ws.Range("N107").Value = -6173.25
# Row 132
$ws.Range("H132").Value = 2338.4827
$ws.Range("I132").Value = 1108.6666
$ws.Range("J132").Value = 3206.5881
$ws.Range("K132").Value = 3325.9998
$ws.Range("L132").Value = 9619.764299999999
$ws.Range("M132").Value = -795.9998000000001
$ws.Range("N132").Value = -14679.7643
